# Conveqs hidden release. New output format: grammar sketch Comments for sentence pairs
#
# Adds a third "comments" column to the sentence-pairs template, one
# comment per existing sentence-pair row, widens that column, and nudges
# a couple of cosmetic workbook bits (selection, a tiny 8pt helper font
# that Mac Excel stamps in for its phonetic-guide bookkeeping) to line
# up with the authored state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "comments" column (C), one row per existing sentence pair ---
$ws.Range("C1").Value = "comments"
$ws.Range("C2").Value = "comments about sentence 1"
$ws.Range("C3").Value = "comments about sentence 2"
$ws.Range("C4").Value = "comments about sentence 3"
$ws.Range("C5").Value = "comments about sentence 4"

# Give the new column its own width (narrower than the source/target ones).
# Excel quantizes ColumnWidth to whole pixels, so this is the closest
# reachable value to the authored 34.1640625 (resolves to ~34.1667).
$ws.Columns.Item(3).ColumnWidth = 33.25

# --- Tiny 8pt font slot (mirrors the phonetic-guide font Mac Excel adds) ---
# Stamped on a scratch sheet so it lands in the shared style table without
# leaving any formatting on real data, then the scratch sheet is discarded.
# (Calibri is already the sheet's default font, so only the size needs to
# change - touching .Name too would mint a second, redundant font entry.)
$excel.DisplayAlerts = $false
$scratch = $wb.Worksheets.Add()
$scratch.Range("A1").Font.Size = 8
$scratch.Delete()
$excel.DisplayAlerts = $true

# Reflect where the author's cursor ended up after the edit
$ws.Range("C10").Select()
